# Update the "p" column (4th column) values in the Mantel correlogram table.
# Row numbers are 1-based and include the header row.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 2;  Old = "0.154"; New = "0.156" },
    @{ Row = 3;  Old = "0.138"; New = "0.144" },
    @{ Row = 4;  Old = "0.248"; New = "0.242" },
    @{ Row = 5;  Old = "0.06";  New = "0.064" },
    @{ Row = 6;  Old = "0.372"; New = "0.363" },
    @{ Row = 7;  Old = "0.496"; New = "0.484" },
    @{ Row = 8;  Old = "0.619"; New = "0.604" },
    @{ Row = 9;  Old = "0.154"; New = "0.128" },
    @{ Row = 10; Old = "0.821"; New = "0.753" },
    @{ Row = 13; Old = "0.71";  New = "0.709" }
)

foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, 4)
    $range = $cell.Range
    # Cell.Range.Text includes the trailing end-of-cell marker (CR + BEL),
    # so trim trailing control characters before comparing.
    $current = $range.Text.TrimEnd([char]13, [char]7)
    if ($current -eq $change.Old) {
        $range.Text = $change.New
    } else {
        Write-Output ("Mismatch at row " + $change.Row + ": expected '" + $change.Old + "' but found '" + $current + "'")
    }
}
